$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.265.24"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").Value = "2.790.79"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'346.55"
$ws.Range("E5").Value = "  +4.14%  "
$ws.Range("D6").Value = "'116.10"
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("E7").Value = "  +3.47%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.589"
$ws.Range("E9").Value = "  +3.14%  "
$ws.Range("D10").Value = "'42.82"
$ws.Range("E10").Value = "  +3.30%  "
$ws.Range("D11").Value = "'0.0858"
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("D12").Value = "'20.02"
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").Value = "'7.84"
$ws.Range("E14").Value = "  +2.53%  "
$ws.Range("D15").Value = "3.228.51"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "2.772.19"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "52.107.63"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").Value = "'3.21"
$ws.Range("E19").Value = "  +6.67%  "
$ws.Range("D20").Value = "'7.21"
$ws.Range("E20").Value = "  +4.91%  "
$ws.Range("D21").Value = "'13.40"
$ws.Range("E21").Value = "  -3.07%  "
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("D23").Value = "'69.97"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "'269.23"
$ws.Range("E24").Value = "  -3.85%  "
$ws.Range("D25").Value = "'2.75"
$ws.Range("E25").Value = "  +4.28%  "
$ws.Range("D26").Value = "'26.63"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D31").Value = "'34.80"
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("D32").Value = "'50.20"
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").Value = "'0.0451"
$ws.Range("E33").Value = "  +28.47%  "
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'2.11"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").Value = "'4.97"
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").Value = "'18.70"
$ws.Range("E39").Value = "  -4.68%  "
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").Value = "'2.66"
$ws.Range("E41").Value = "  +16.46%  "
$ws.Range("D42").Value = "'128.02"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("D43").Value = "'23.28"
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("E46").Value = "  -2.23%  "
$ws.Range("D47").Value = "2.065.79"
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("E48").Value = "  +2.64%  "
$ws.Range("D49").Value = "'0.966"
$ws.Range("E49").Value = "  +12.65%  "
$ws.Range("D50").Value = "'5.52"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").Value = "'8.90"
$ws.Range("E51").Value = "  -1.64%  "
